$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add three new rows of UI work log entries
$ws.Range("E18").Value = [DateTime]"2019-07-25"
$ws.Range("F18").Value = 0.5
$ws.Range("G18").Value = "Anregungen für die Erstellung einer UI gesucht"

$ws.Range("E19").Value = [DateTime]"2019-07-28"
$ws.Range("F19").Value = 1.5
$ws.Range("G19").Value = "Erstellung einer ersten UI mit SceneBuilder"

$ws.Range("E20").Value = [DateTime]"2019-07-29"
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = "Weitere UI-Designs mit SceneBuilder erstellt, insgesamt 3 Designs fertiggestellt, Statusupdate gegeben"

# Copy style from row 17 to the new rows 18-20 (date/hours/text formatting)
$ws.Range("E17:G17").Copy()
$ws.Range("E18:G20").PasteSpecial(-4122)

# Update selection to match the recorded state
$ws.Range("G5").Select()
